$wb = $excel.ActiveWorkbook

# Sheet "ip_address_list" (sheet1): update the notes cell D1 and add D2
$wsIp = $wb.Worksheets.Item("ip_address_list")
$wsIp.Range("D1").Value = "poznggv`nf`nfhk`nhvj`nhg`nguk`nf`nf`nf`nf`njjjf`nf"
$wsIp.Range("D2").Value = "f`nf`nf`nf"

# Sheet "disk_list" (sheet3): extend F4 note with extra lines
$wsDisk = $wb.Worksheets.Item("disk_list")
$wsDisk.Range("F4").Value = "druhá síť Valeo`nfg`nf`nf"

# Sheet "Settings" (sheet4): flip a few default-behavior flags (grid -> pack)
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B1").Value = 0
$wsSettings.Range("B3").Value = 1
$wsSettings.Range("B4").Value = 1
